$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Each of the nine files is structured as follows:"
#    -> "Each of the nine files contains the following sheets:"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Each of the nine files is structured as follows:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Each of the nine files contains the following sheets:", 2
) | Out-Null

# ---------------------------------------------------------------------------
# 2) STATS INDICATORS sheet description: the CCI paragraph block already
#    ends with "...columns from K to P". Add a parallel block describing
#    the RSI indicator columns (Q to AK), mirroring the CCI block just
#    above it.
# ---------------------------------------------------------------------------

# Locate the anchor paragraph - the last paragraph in the document:
# "The same applies for the 20-days-period setup in columns from K to P"
$anchor = $d.Paragraphs($d.Paragraphs.Count)

# Helper: italicize a given literal substring inside a paragraph's range.
function Set-ItalicSubstring {
    param($paragraph, [string]$needle)
    $searchRange = $paragraph.Range.Duplicate()
    $searchRange.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $searchRange.Italic = 1
}

# --- new paragraph: "Columns from Q to AK refer to the RSI indicator." (ilvl 0) ---
$anchor.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs($anchor.Index + 1)
$p1.Range.Text = "Columns from Q to AK refer to the RSI indicator.`r"
$p1.Range.ListFormat.ListLevelNumber = 1

# --- new paragraph: "Column Q (T) represents the average value for the overbought (oversold) status in the 7 days-period setup." (ilvl 1) ---
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs($p1.Index + 1)
$p2.Range.Text = "Column Q (T) represents the average value for the overbought (oversold) status in the 7 days-period setup.`r"
$p2.Range.ListFormat.ListLevelNumber = 2
Set-ItalicSubstring $p2 "T"
Set-ItalicSubstring $p2 "oversold"

# --- new paragraph: "Column R (U) represents the maximum value for the overbought (oversold) status in the 7 days-period setup." (ilvl 1) ---
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs($p2.Index + 1)
$p3.Range.Text = "Column R (U) represents the maximum value for the overbought (oversold) status in the 7 days-period setup.`r"
$p3.Range.ListFormat.ListLevelNumber = 2
Set-ItalicSubstring $p3 "U"
Set-ItalicSubstring $p3 "oversold"

# --- new paragraph: "Column S (V) represents the minimum value for the overbought (oversold) status in the 7 days-period setup." (ilvl 1) ---
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs($p3.Index + 1)
$p4.Range.Text = "Column S (V) represents the minimum value for the overbought (oversold) status in the 7 days-period setup.`r"
$p4.Range.ListFormat.ListLevelNumber = 2
Set-ItalicSubstring $p4 "V"
Set-ItalicSubstring $p4 "oversold"

# --- new paragraph: "Column W represents the number of crossings of the 50-line in the 7 days-period setup." (ilvl 1) ---
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs($p4.Index + 1)
$p5.Range.Text = "Column W represents the number of crossings of the 50-line in the 7 days-period setup.`r"
$p5.Range.ListFormat.ListLevelNumber = 2

# --- new paragraph: "The same applies for the 14-days-period setup in columns from X to AD." (ilvl 1) ---
$p5.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs($p5.Index + 1)
$p6.Range.Text = "The same applies for the 14-days-period setup in columns from X to AD.`r"
$p6.Range.ListFormat.ListLevelNumber = 2

# --- new paragraph: "The same applies for the 21-days-period setup in columns from AE to AK." (ilvl 1, last paragraph) ---
$p6.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs($p6.Index + 1)
$p7.Range.Text = "The same applies for the 21-days-period setup in columns from AE to AK.`r"
$p7.Range.ListFormat.ListLevelNumber = 2
